# "small refactor & docs"
#
# Journal de travail: log a new day of work (2023-07-21) in the first
# still-empty row (row 86) of the "Journal de travail" table, and tidy up
# the stray formatting that was sitting on the following empty rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Journal de travail")

# Rows 87-92 only carried leftover number-format styling on column C with
# no data; clear that formatting entirely (row 86 will get real content
# below, rows 87-92 stay blank/unstyled).
$ws.Range("C86:C92").Clear()

# New journal entry for 2023-07-21: 9h spent on "Rédaction" (Rapport).
$ws.Range("A86").Value = "07/21/2023"
$ws.Range("B86").Value = "Rédaction"
$ws.Range("C86").Value = 9
$ws.Range("D86").Value = "Rapport"

# Leave the selection where it was in the saved workbook (next empty row).
$ws.Range("D87").Select() | Out-Null
